$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking strings in the D:E data range
# so Excel does not auto-convert them to numbers (which would change
# both the stored type and introduce floating-point rounding).
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.270.14"
$ws.Range("E2").Value = "  +2.70%  "

# Row 3
$ws.Range("D3").Value = "1.894.46"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.82%  "

# Row 5
$ws.Range("D5").Value = "315.05"
$ws.Range("E5").Value = "  -0.51%  "

# Row 6
$ws.Range("E6").Value = "  -0.92%  "

# Row 7
$ws.Range("D7").Value = "0.5142"
$ws.Range("E7").Value = "  +0.48%  "

# Row 8
$ws.Range("D8").Value = "0.3925"
$ws.Range("E8").Value = "  -0.96%  "

# Row 9
$ws.Range("D9").Value = "0.08422"
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").Value = "42.47"
$ws.Range("E10").Value = "  +1.61%  "

# Row 11
$ws.Range("D11").Value = "1.114"
$ws.Range("E11").Value = "  +0.33%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.899.31"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "6.242"
$ws.Range("E13").Value = "  -0.46%  "

# Row 14
$ws.Range("D14").Value = "20.67"
$ws.Range("E14").Value = "  +0.78%  "

# Row 15
$ws.Range("D15").Value = "7.313"
$ws.Range("E15").Value = "  +0.43%  "

# Row 16
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.87%  "

# Row 17
$ws.Range("D17").Value = "93.10"
$ws.Range("E17").Value = "  +1.99%  "

# Row 18
$ws.Range("D18").Value = "0.00001107"
$ws.Range("E18").Value = "  -0.22%  "

# Row 19
$ws.Range("D19").Value = "0.06751"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("D20").Value = "17.84"
$ws.Range("E20").Value = "  +0.61%  "

# Row 21
$ws.Range("E21").Value = "  -0.85%  "

# Row 22
$ws.Range("D22").Value = "6.015"
$ws.Range("E22").Value = "  +1.05%  "

# Row 23
$ws.Range("D23").Value = "29.290.31"
$ws.Range("E23").Value = "  +2.60%  "

# Row 24
$ws.Range("D24").Value = "11.14"
$ws.Range("E24").Value = "  -0.06%  "

# Row 25
$ws.Range("E25").Value = "  -2.58%  "

# Row 26
$ws.Range("D26").Value = "2.113.85"
$ws.Range("E26").Value = "  +0.23%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "20.98"
$ws.Range("E27").Value = "  +0.63%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "159.15"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29
$ws.Range("D29").Value = "2.434"
$ws.Range("E29").Value = "  +1.97%  "

# Row 30
$ws.Range("D30").Value = "128.04"
$ws.Range("E30").Value = "  +0.59%  "

# Row 31
$ws.Range("D31").Value = "1.058"
$ws.Range("E31").Value = "  +0.94%  "

# Row 32
$ws.Range("E32").Value = "  -0.67%  "

# Row 33
$ws.Range("D33").Value = "6.130"
$ws.Range("E33").Value = "  +6.00%  "

# Row 34
$ws.Range("D34").Value = "3.647"
$ws.Range("E34").Value = "  +0.59%  "

# Row 35
$ws.Range("D35").Value = "0.02476"
$ws.Range("E35").Value = "  +1.40%  "

# Row 36
$ws.Range("D36").Value = "0.06538"
$ws.Range("E36").Value = "  +0.50%  "

# Row 37
$ws.Range("D37").Value = "9.042"
$ws.Range("E37").Value = "  +1.57%  "

# Row 38
$ws.Range("E38").Value = "  +0.53%  "

# Row 39
$ws.Range("D39").Value = "1.228"
$ws.Range("E39").Value = "  +3.15%  "

# Row 40
$ws.Range("D40").Value = "5.128"
$ws.Range("E40").Value = "  +1.37%  "

# Row 41
$ws.Range("D41").Value = "0.6511"
$ws.Range("E41").Value = "  +0.75%  "

# Row 42
$ws.Range("D42").Value = "1.233"
$ws.Range("E42").Value = "  -2.66%  "

# Row 43
$ws.Range("D43").Value = "11.26"
$ws.Range("E43").Value = "  +0.51%  "

# Row 44
$ws.Range("D44").Value = "0.6059"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45
$ws.Range("D45").Value = "13.22"
$ws.Range("E45").Value = "  +1.27%  "

# Row 46
$ws.Range("D46").Value = "3.676"
$ws.Range("E46").Value = "  -1.01%  "

# Row 47
$ws.Range("D47").Value = "2.046"
$ws.Range("E47").Value = "  +1.21%  "

# Row 48
$ws.Range("E48").Value = "  +1.48%  "

# Row 49
$ws.Range("D49").Value = "123.35"
$ws.Range("E49").Value = "  +0.73%  "

# Row 50
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -2.44%  "

# Row 51
$ws.Range("D51").Value = "77.67"
$ws.Range("E51").Value = "  +0.61%  "

# Restore original (default/general) formatting on the data range so the
# saved style indices match the source workbook exactly.
$dataRange.ClearFormats()
